$wb = $excel.ActiveWorkbook

# --- three_line ---
$ws = $wb.Worksheets.Item("three_line")
$ws.Range("A362:A364").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E362:E364").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("G362:G364").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("I362:I364").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A362").Value = 45442.55208333334
$ws.Range("B362").Value = "18-06-2024 09:15:00"
$ws.Range("C362").Value = "hour"
$ws.Range("D362").Value = "TATAINVEST.NS"
$ws.Range("E362").Value = 45400.38541666666
$ws.Range("F362").Value = 7000
$ws.Range("G362").Value = 45422.38541666666
$ws.Range("H362").Value = 6890
$ws.Range("I362").Value = 45436.38541666666
$ws.Range("J362").Value = 6874
$ws.Range("K362").Value = "High"
$ws.Range("L362").Value = "18/06/2024 04:46:11"

$ws.Range("A363").Value = 45450.55208333334
$ws.Range("B363").Value = "18-06-2024 09:15:00"
$ws.Range("C363").Value = "hour"
$ws.Range("D363").Value = "SHILCTECH.BO"
$ws.Range("E363").Value = 45447.38541666666
$ws.Range("F363").Value = 5350
$ws.Range("G363").Value = 45447.42708333334
$ws.Range("H363").Value = 5350
$ws.Range("I363").Value = 45449.38541666666
$ws.Range("J363").Value = 5394
$ws.Range("K363").Value = "High"
$ws.Range("L363").Value = "18/06/2024 04:46:11"

$ws.Range("A364").Value = 45440.55208333334
$ws.Range("B364").Value = "18-06-2024 09:15:00"
$ws.Range("C364").Value = "hour"
$ws.Range("D364").Value = "STOVEKRAFT.NS"
$ws.Range("E364").Value = 45412.46875
$ws.Range("F364").Value = 464.2999877929688
$ws.Range("G364").Value = 45433.55208333334
$ws.Range("H364").Value = 509.8500061035156
$ws.Range("I364").Value = 45439.38541666666
$ws.Range("J364").Value = 517
$ws.Range("K364").Value = "High"
$ws.Range("L364").Value = "18/06/2024 04:46:11"

# --- two_line ---
$ws = $wb.Worksheets.Item("two_line")
$ws.Range("A185:A190").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E185:E190").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("G185:G190").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A185").Value = 45406.55208333334
$ws.Range("B185").Value = "18-06-2024 10:15:00"
$ws.Range("C185").Value = "hour"
$ws.Range("D185").Value = "DMART.NS"
$ws.Range("E185").Value = 45392.38541666666
$ws.Range("F185").Value = 4831.5
$ws.Range("G185").Value = 45401.38541666666
$ws.Range("H185").Value = 4831.7998046875
$ws.Range("I185").Value = "High"
$ws.Range("J185").Value = "18/06/2024 04:46:11"

$ws.Range("A186").Value = 45455.63541666666
$ws.Range("B186").Value = "18-06-2024 10:15:00"
$ws.Range("C186").Value = "hour"
$ws.Range("D186").Value = "NAUKRI.NS"
$ws.Range("E186").Value = 45450.38541666666
$ws.Range("F186").Value = 6359.9501953125
$ws.Range("G186").Value = 45454.46875
$ws.Range("H186").Value = 6327.60009765625
$ws.Range("I186").Value = "High"
$ws.Range("J186").Value = "18/06/2024 04:46:11"

$ws.Range("A187").Value = 45450.55208333334
$ws.Range("B187").Value = "18-06-2024 09:15:00"
$ws.Range("C187").Value = "hour"
$ws.Range("D187").Value = "SHILCTECH.BO"
$ws.Range("E187").Value = 45447.38541666666
$ws.Range("F187").Value = 5350
$ws.Range("G187").Value = 45449.38541666666
$ws.Range("H187").Value = 5394
$ws.Range("I187").Value = "High"
$ws.Range("J187").Value = "18/06/2024 04:46:11"

$ws.Range("A188").Value = 45450.55208333334
$ws.Range("B188").Value = "18-06-2024 09:15:00"
$ws.Range("C188").Value = "hour"
$ws.Range("D188").Value = "SHILCTECH.BO"
$ws.Range("E188").Value = 45447.42708333334
$ws.Range("F188").Value = 5350
$ws.Range("G188").Value = 45449.38541666666
$ws.Range("H188").Value = 5394
$ws.Range("I188").Value = "High"
$ws.Range("J188").Value = "18/06/2024 04:46:11"

$ws.Range("A189").Value = 45427.42708333334
$ws.Range("B189").Value = "18-06-2024 09:15:00"
$ws.Range("C189").Value = "hour"
$ws.Range("D189").Value = "KPEL.BO"
$ws.Range("E189").Value = 45425.46875
$ws.Range("F189").Value = 422.1000061035156
$ws.Range("G189").Value = 45425.55208333334
$ws.Range("H189").Value = 422.1000061035156
$ws.Range("I189").Value = "Low"
$ws.Range("J189").Value = "18/06/2024 04:46:11"

$ws.Range("A190").Value = 45439.42708333334
$ws.Range("B190").Value = "18-06-2024 09:15:00"
$ws.Range("C190").Value = "hour"
$ws.Range("D190").Value = "CARYSIL.NS"
$ws.Range("E190").Value = 45421.59375
$ws.Range("F190").Value = 847.0499877929688
$ws.Range("G190").Value = 45435.55208333334
$ws.Range("H190").Value = 850
$ws.Range("I190").Value = "Low"
$ws.Range("J190").Value = "18/06/2024 04:46:11"

# --- ph_pl_breakout_line ---
$ws = $wb.Worksheets.Item("ph_pl_breakout_line")
$ws.Range("B1299:B1316").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A1299").Value = "BAJAJST.BO"
$ws.Range("B1299").Value = 45446.38541666666
$ws.Range("C1299").Value = 1500
$ws.Range("D1299").Value = 1420
$ws.Range("E1299").Value = 1467.75
$ws.Range("F1299").Value = "High"
$ws.Range("G1299").Value = 1500
$ws.Range("H1299").Value = "hour"
$ws.Range("I1299").Value = "18-06-2024 09:15:00"
$ws.Range("J1299").Value = 1510
$ws.Range("K1299").Value = 1495
$ws.Range("L1299").Value = "18/06/2024 04:46:11"

$ws.Range("A1300").Value = "TRF.NS"
$ws.Range("B1300").Value = 45456.38541666666
$ws.Range("C1300").Value = 648.9500122070312
$ws.Range("D1300").Value = 611.25
$ws.Range("E1300").Value = 630
$ws.Range("F1300").Value = "Low"
$ws.Range("G1300").Value = 611.25
$ws.Range("H1300").Value = "hour"
$ws.Range("I1300").Value = "18-06-2024 09:15:00"
$ws.Range("J1300").Value = 610.8499755859375
$ws.Range("K1300").Value = 623.2999877929688
$ws.Range("L1300").Value = "18/06/2024 04:46:11"

$ws.Range("A1301").Value = "ASMTEC.BO"
$ws.Range("B1301").Value = 45455.38541666666
$ws.Range("C1301").Value = 1149.949951171875
$ws.Range("D1301").Value = 1101.099975585938
$ws.Range("E1301").Value = 1129.849975585938
$ws.Range("F1301").Value = "High"
$ws.Range("G1301").Value = 1149.949951171875
$ws.Range("H1301").Value = "hour"
$ws.Range("I1301").Value = "18-06-2024 09:15:00"
$ws.Range("J1301").Value = 1160
$ws.Range("K1301").Value = 1122
$ws.Range("L1301").Value = "18/06/2024 04:46:11"

$ws.Range("A1302").Value = "SPELS.BO"
$ws.Range("B1302").Value = 45455.38541666666
$ws.Range("C1302").Value = 143.9499969482422
$ws.Range("D1302").Value = 136
$ws.Range("E1302").Value = 139.5
$ws.Range("F1302").Value = "High"
$ws.Range("G1302").Value = 143.9499969482422
$ws.Range("H1302").Value = "hour"
$ws.Range("I1302").Value = "18-06-2024 09:15:00"
$ws.Range("J1302").Value = 145.5
$ws.Range("K1302").Value = 143
$ws.Range("L1302").Value = "18/06/2024 04:46:11"

$ws.Range("A1303").Value = "SPELS.BO"
$ws.Range("B1303").Value = 45455.46875
$ws.Range("C1303").Value = 138.8000030517578
$ws.Range("D1303").Value = 135
$ws.Range("E1303").Value = 137
$ws.Range("F1303").Value = "Low"
$ws.Range("G1303").Value = 135
$ws.Range("H1303").Value = "hour"
$ws.Range("I1303").Value = "18-06-2024 09:15:00"
$ws.Range("J1303").Value = 134.75
$ws.Range("K1303").Value = 140
$ws.Range("L1303").Value = "18/06/2024 04:46:11"

$ws.Range("A1304").Value = "SPELS.BO"
$ws.Range("B1304").Value = 45455.51041666666
$ws.Range("C1304").Value = 139
$ws.Range("D1304").Value = 135
$ws.Range("E1304").Value = 136.8999938964844
$ws.Range("F1304").Value = "Low"
$ws.Range("G1304").Value = 135
$ws.Range("H1304").Value = "hour"
$ws.Range("I1304").Value = "18-06-2024 09:15:00"
$ws.Range("J1304").Value = 134.75
$ws.Range("K1304").Value = 140
$ws.Range("L1304").Value = "18/06/2024 04:46:11"

$ws.Range("A1305").Value = "CROWN.NS"
$ws.Range("B1305").Value = 45421.38541666666
$ws.Range("C1305").Value = 260.4500122070312
$ws.Range("D1305").Value = 236.75
$ws.Range("E1305").Value = 236.75
$ws.Range("F1305").Value = "High"
$ws.Range("G1305").Value = 260.4500122070312
$ws.Range("H1305").Value = "hour"
$ws.Range("I1305").Value = "18-06-2024 09:15:00"
$ws.Range("J1305").Value = 262.4800109863281
$ws.Range("K1305").Value = 249.9900054931641
$ws.Range("L1305").Value = "18/06/2024 04:46:11"

$ws.Range("A1306").Value = "ZODIAC.NS"
$ws.Range("B1306").Value = 45448.55208333334
$ws.Range("C1306").Value = 611.8499755859375
$ws.Range("D1306").Value = 611.8499755859375
$ws.Range("E1306").Value = 611.8499755859375
$ws.Range("F1306").Value = "Low"
$ws.Range("G1306").Value = 611.8499755859375
$ws.Range("H1306").Value = "hour"
$ws.Range("I1306").Value = "18-06-2024 09:15:00"
$ws.Range("J1306").Value = 608.9500122070312
$ws.Range("K1306").Value = 621.3499755859375
$ws.Range("L1306").Value = "18/06/2024 04:46:11"

$ws.Range("A1307").Value = "ZODIAC.NS"
$ws.Range("B1307").Value = 45448.59375
$ws.Range("C1307").Value = 611.8499755859375
$ws.Range("D1307").Value = 611.8499755859375
$ws.Range("E1307").Value = 611.8499755859375
$ws.Range("F1307").Value = "Low"
$ws.Range("G1307").Value = 611.8499755859375
$ws.Range("H1307").Value = "hour"
$ws.Range("I1307").Value = "18-06-2024 09:15:00"
$ws.Range("J1307").Value = 608.9500122070312
$ws.Range("K1307").Value = 621.3499755859375
$ws.Range("L1307").Value = "18/06/2024 04:46:11"

$ws.Range("A1308").Value = "ZODIAC.NS"
$ws.Range("B1308").Value = 45456.38541666666
$ws.Range("C1308").Value = 611.0999755859375
$ws.Range("D1308").Value = 611.0999755859375
$ws.Range("E1308").Value = 611.0999755859375
$ws.Range("F1308").Value = "Low"
$ws.Range("G1308").Value = 611.0999755859375
$ws.Range("H1308").Value = "hour"
$ws.Range("I1308").Value = "18-06-2024 09:15:00"
$ws.Range("J1308").Value = 608.9500122070312
$ws.Range("K1308").Value = 621.3499755859375
$ws.Range("L1308").Value = "18/06/2024 04:46:11"

$ws.Range("A1309").Value = "CENTRALBK.BO"
$ws.Range("B1309").Value = 45453.38541666666
$ws.Range("C1309").Value = 65.7300033569336
$ws.Range("D1309").Value = 64.75
$ws.Range("E1309").Value = 65.30000305175781
$ws.Range("F1309").Value = "High"
$ws.Range("G1309").Value = 65.7300033569336
$ws.Range("H1309").Value = "hour"
$ws.Range("I1309").Value = "18-06-2024 09:15:00"
$ws.Range("J1309").Value = 65.80000305175781
$ws.Range("K1309").Value = 65.55000305175781
$ws.Range("L1309").Value = "18/06/2024 04:46:11"

$ws.Range("A1310").Value = "SIL.NS"
$ws.Range("B1310").Value = 45446.38541666666
$ws.Range("C1310").Value = 24.89999961853027
$ws.Range("D1310").Value = 23.45000076293945
$ws.Range("E1310").Value = 23.85000038146973
$ws.Range("F1310").Value = "High"
$ws.Range("G1310").Value = 24.89999961853027
$ws.Range("H1310").Value = "hour"
$ws.Range("I1310").Value = "18-06-2024 09:15:00"
$ws.Range("J1310").Value = 24.94000053405762
$ws.Range("K1310").Value = 24.85000038146973
$ws.Range("L1310").Value = "18/06/2024 04:46:11"

$ws.Range("A1311").Value = "AGARIND.NS"
$ws.Range("B1311").Value = 45455.42708333334
$ws.Range("C1311").Value = 1143
$ws.Range("D1311").Value = 1077
$ws.Range("E1311").Value = 1120
$ws.Range("F1311").Value = "High"
$ws.Range("G1311").Value = 1143
$ws.Range("H1311").Value = "hour"
$ws.Range("I1311").Value = "18-06-2024 09:15:00"
$ws.Range("J1311").Value = 1179.300048828125
$ws.Range("K1311").Value = 1142
$ws.Range("L1311").Value = "18/06/2024 04:46:11"

$ws.Range("A1312").Value = "AGARIND.NS"
$ws.Range("B1312").Value = 45455.46875
$ws.Range("C1312").Value = 1143
$ws.Range("D1312").Value = 1111.449951171875
$ws.Range("E1312").Value = 1116.900024414062
$ws.Range("F1312").Value = "High"
$ws.Range("G1312").Value = 1143
$ws.Range("H1312").Value = "hour"
$ws.Range("I1312").Value = "18-06-2024 09:15:00"
$ws.Range("J1312").Value = 1179.300048828125
$ws.Range("K1312").Value = 1142
$ws.Range("L1312").Value = "18/06/2024 04:46:11"

$ws.Range("A1313").Value = "RPEL.BO"
$ws.Range("B1313").Value = 45439.59375
$ws.Range("C1313").Value = 725.4500122070312
$ws.Range("D1313").Value = 682.4000244140625
$ws.Range("E1313").Value = 710
$ws.Range("F1313").Value = "High"
$ws.Range("G1313").Value = 725.4500122070312
$ws.Range("H1313").Value = "hour"
$ws.Range("I1313").Value = "18-06-2024 09:15:00"
$ws.Range("J1313").Value = 729
$ws.Range("K1313").Value = 718.9000244140625
$ws.Range("L1313").Value = "18/06/2024 04:46:11"

$ws.Range("A1314").Value = "SANJIVIN.BO"
$ws.Range("B1314").Value = 45454.38541666666
$ws.Range("C1314").Value = 176.8000030517578
$ws.Range("D1314").Value = 170.1999969482422
$ws.Range("E1314").Value = 170.5500030517578
$ws.Range("F1314").Value = "High"
$ws.Range("G1314").Value = 176.8000030517578
$ws.Range("H1314").Value = "hour"
$ws.Range("I1314").Value = "18-06-2024 09:15:00"
$ws.Range("J1314").Value = 177.8999938964844
$ws.Range("K1314").Value = 175
$ws.Range("L1314").Value = "18/06/2024 04:46:11"

$ws.Range("A1315").Value = "UNIAUTO.BO"
$ws.Range("B1315").Value = 45443.38541666666
$ws.Range("C1315").Value = 171.9499969482422
$ws.Range("D1315").Value = 164
$ws.Range("E1315").Value = 167.4499969482422
$ws.Range("F1315").Value = "High"
$ws.Range("G1315").Value = 171.9499969482422
$ws.Range("H1315").Value = "hour"
$ws.Range("I1315").Value = "18-06-2024 09:15:00"
$ws.Range("J1315").Value = 174
$ws.Range("K1315").Value = 170.8500061035156
$ws.Range("L1315").Value = "18/06/2024 04:46:11"

$ws.Range("A1316").Value = "UNIAUTO.BO"
$ws.Range("B1316").Value = 45454.38541666666
$ws.Range("C1316").Value = 171.9499969482422
$ws.Range("D1316").Value = 163.1999969482422
$ws.Range("E1316").Value = 166.5
$ws.Range("F1316").Value = "High"
$ws.Range("G1316").Value = 171.9499969482422
$ws.Range("H1316").Value = "hour"
$ws.Range("I1316").Value = "18-06-2024 09:15:00"
$ws.Range("J1316").Value = 174
$ws.Range("K1316").Value = 170.8500061035156
$ws.Range("L1316").Value = "18/06/2024 04:46:11"

